$wb = $excel.ActiveWorkbook

$msg = "Validado com sucesso! Nenhuma divergência entre o SPED e o relatório foi encontrada!"

# --- Sheet 1: Bico ---
$ws1 = $wb.Worksheets.Item("Bico")
$ws1.Range("H1").Value = "Obs_relatorio"
$ws1.Range("I1").Value = "Obs_sped"
for ($r = 2; $r -le 7; $r++) {
    $ws1.Cells.Item($r, 8).Value = $msg
    $ws1.Cells.Item($r, 9).Value = ""
}

# --- Sheet 2: Tanque ---
$ws2 = $wb.Worksheets.Item("Tanque")
$ws2.Range("F1").Value = "Obs_relatorio"
$ws2.Range("G1").Value = "Obs_sped"
for ($r = 2; $r -le 4; $r++) {
    $ws2.Cells.Item($r, 6).Value = $msg
    $ws2.Cells.Item($r, 7).Value = ""
}
